# Efectividades.xlsx - "Finalizada tabla de efectividades"
# Fills in the remaining pairwise subject-comparison cells (Sujeto 6..9
# combinations) that were still blank, and updates the view/selection
# state to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sujeto 6 block (row 17-19): vs Sujeto 7 / Sujeto 8 / Sujeto 9 ---
$ws.Range("AG17").Value = 1
$ws.Range("AJ17").Value = 1
$ws.Range("AM17").Value = 0.98333300000000001
$ws.Range("AP17").Value = 1
$ws.Range("AS17").Value = 0.97916599999999998
$ws.Range("AV17").Value = 1

$ws.Range("AF19").Value = 39
$ws.Range("AG19").Value = 56
$ws.Range("AH19").Value = 33
$ws.Range("AI19").Value = 24
$ws.Range("AJ19").Value = 49
$ws.Range("AL19").Value = 38
$ws.Range("AM19").Value = 49
$ws.Range("AO19").Value = 38
$ws.Range("AP19").Value = 49
$ws.Range("AR19").Value = 42
$ws.Range("AS19").Value = 49
$ws.Range("AU19").Value = 46
$ws.Range("AV19").Value = 49

# --- Sujeto 7 block (row 20-22): vs Sujeto 8 / Sujeto 9 ---
$ws.Range("AM20").Value = 1
$ws.Range("AP20").Value = 1
$ws.Range("AS20").Value = 0.98958299999999999

$ws.Range("AL22").Value = 30
$ws.Range("AM22").Value = 49
$ws.Range("AO22").Value = 38
$ws.Range("AP22").Value = 49
$ws.Range("AR22").Value = 31
$ws.Range("AS22").Value = 38

# --- Sujeto 8 block (row 23-25): vs Sujeto 9 ---
$ws.Range("AS23").Value = 0.97058800000000001
$ws.Range("AV23").Value = 0.94117600000000001

$ws.Range("AR25").Value = 41
$ws.Range("AS25").Value = 49
$ws.Range("AU25").Value = 40
$ws.Range("AV25").Value = 49

# --- View state: scroll back to the left edge, zoom to 87%, and leave
# the selection on A27 (the effectiveness-summary row) ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 87
$ws.Range("A27").Select()
